{"js": "// Replace each \"old\" division-problem answer string with its \"new\" value.\n// Every old string in this table is unique within the document, so a\n// simple search-and-replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"152\u00f72=76, 0\", \"159\u00f77=22, 5\"],\n  [\"303\u00f72=151, 1\", \"109\u00f74=27, 1\"],\n  [\"485\u00f78=60, 5\", \"361\u00f72=180, 1\"],\n  [\"752\u00f78=94, 0\", \"534\u00f79=59, 3\"],\n  [\"429\u00f79=47, 6\", \"371\u00f75=74, 1\"],\n  [\"170\u00f77=24, 2\", \"853\u00f72=426, 1\"],\n  [\"388\u00f77=55, 3\", \"265\u00f74=66, 1\"],\n  [\"291\u00f72=145, 1\", \"246\u00f76=41, 0\"],\n  [\"869\u00f74=217, 1\", \"134\u00f75=26, 4\"],\n  [\"494\u00f72=247, 0\", \"912\u00f79=101, 3\"],\n  [\"586\u00f76=97, 4\", \"947\u00f73=315, 2\"],\n  [\"776\u00f76=129, 2\", \"331\u00f74=82, 3\"],\n  [\"612\u00f73=204, 0\", \"414\u00f79=46, 0\"],\n  [\"711\u00f79=79, 0\", \"364\u00f77=52, 0\"],\n  [\"318\u00f74=79, 2\", \"940\u00f72=470, 0\"],\n  [\"155\u00f76=25, 5\", \"816\u00f75=163, 1\"],\n  [\"135\u00f78=16, 7\", \"534\u00f72=267, 0\"],\n  [\"848\u00f78=106, 0\", \"873\u00f73=291, 0\"],\n  [\"762\u00f78=95, 2\", \"899\u00f78=112, 3\"],\n  [\"145\u00f72=72, 1\", \"400\u00f77=57, 1\"],\n  [\"281\u00f77=40, 1\", \"271\u00f74=67, 3\"],\n  [\"757\u00f75=151, 2\", \"547\u00f79=60, 7\"],\n  [\"512\u00f75=102, 2\", \"581\u00f78=72, 5\"],\n  [\"995\u00f75=199, 0\", \"287\u00f75=57, 2\"],\n  [\"407\u00f79=45, 2\", \"613\u00f73=204, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"old\" division-problem answer string with its \"new\" value.\n# Every old string below is unique within the document, so a simple\n# Find/Replace per pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"152\u00f72=76, 0\", \"159\u00f77=22, 5\"),\n    @(\"303\u00f72=151, 1\", \"109\u00f74=27, 1\"),\n    @(\"485\u00f78=60, 5\", \"361\u00f72=180, 1\"),\n    @(\"752\u00f78=94, 0\", \"534\u00f79=59, 3\"),\n    @(\"429\u00f79=47, 6\", \"371\u00f75=74, 1\"),\n    @(\"170\u00f77=24, 2\", \"853\u00f72=426, 1\"),\n    @(\"388\u00f77=55, 3\", \"265\u00f74=66, 1\"),\n    @(\"291\u00f72=145, 1\", \"246\u00f76=41, 0\"),\n    @(\"869\u00f74=217, 1\", \"134\u00f75=26, 4\"),\n    @(\"494\u00f72=247, 0\", \"912\u00f79=101, 3\"),\n    @(\"586\u00f76=97, 4\", \"947\u00f73=315, 2\"),\n    @(\"776\u00f76=129, 2\", \"331\u00f74=82, 3\"),\n    @(\"612\u00f73=204, 0\", \"414\u00f79=46, 0\"),\n    @(\"711\u00f79=79, 0\", \"364\u00f77=52, 0\"),\n    @(\"318\u00f74=79, 2\", \"940\u00f72=470, 0\"),\n    @(\"155\u00f76=25, 5\", \"816\u00f75=163, 1\"),\n    @(\"135\u00f78=16, 7\", \"534\u00f72=267, 0\"),\n    @(\"848\u00f78=106, 0\", \"873\u00f73=291, 0\"),\n    @(\"762\u00f78=95, 2\", \"899\u00f78=112, 3\"),\n    @(\"145\u00f72=72, 1\", \"400\u00f77=57, 1\"),\n    @(\"281\u00f77=40, 1\", \"271\u00f74=67, 3\"),\n    @(\"757\u00f75=151, 2\", \"547\u00f79=60, 7\"),\n    @(\"512\u00f75=102, 2\", \"581\u00f78=72, 5\"),\n    @(\"995\u00f75=199, 0\", \"287\u00f75=57, 2\"),\n    @(\"407\u00f79=45, 2\", \"613\u00f73=204, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
